# Update the StructureDefinition match-criteria workbook:
#  - Metadata sheet: URL, Version, Date, Publisher
#  - Elements sheet: clear stale Constraint(s) on the Extension row,
#    and refresh the Fixed Value URL on the Extension.url row.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-criteria"
$metadata.Range("B3").Value = "8.0.0"
$metadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$metadata.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-criteria"
